# Apply the edits described by the commit:
#  - Update the test-data username on the "LogIn" sheet
#  - Move the active/selected tab from "Register" to "LogIn"
#  - Update the remembered selection on each sheet

$wb = $excel.ActiveWorkbook

$wsRegister = $wb.Worksheets.Item("Register")
$wsLogin    = $wb.Worksheets.Item("LogIn")

# Update the shared test-data value everywhere it is used (Register!I2 and
# LogIn!A2 both held the same "testSelenium29409" string).
$wsRegister.Range("I2").Value = "testDemo1990092"
$wsLogin.Range("A2").Value = "testDemo1990092"

# Register keeps a remembered selection of C22, but is no longer the active tab
$wsRegister.Range("C22").Select()

# LogIn becomes the active sheet/tab, with a remembered selection of A6
$wsLogin.Select()
$wsLogin.Range("A6").Select()
